$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the LBNDIND -> LBNRIND typo in the header row (F1)
$ws.Range("F1").Value = "LBNRIND"

# Update selection to M9 (no data there), matching the post-edit click
$ws.Range("M9").Select()
